$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.291.09'
$ws.Range("E2").Value = '  +2.32%  '

$ws.Range("D3").Value = '2.534.04'
$ws.Range("E3").Value = '  +2.87%  '

$ws.Range("E4").Value = '  +0.07%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '582.84'
$c.ClearFormats()
$ws.Range("E5").Value = '  +1.33%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '153.05'
$c.ClearFormats()
$ws.Range("E6").Value = '  +4.85%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  +0.40%  '

$ws.Range("E9").Value = '  +1.23%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.162'
$c.ClearFormats()
$ws.Range("E10").Value = '  -0.09%  '

$ws.Range("E11").Value = '  +0.33%  '

$ws.Range("E12").Value = '  +0.16%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '29.70'
$c.ClearFormats()
$ws.Range("E13").Value = '  +2.49%  '

$ws.Range("E14").Value = '  +1.48%  '

$ws.Range("D15").Value = '2.994.57'
$ws.Range("E15").Value = '  +2.91%  '

$ws.Range("D16").Value = '63.640.05'
$ws.Range("E16").Value = '  +1.33%  '

$ws.Range("D17").Value = '2.547.94'
$ws.Range("E17").Value = '  +3.32%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '7.91'
$c.ClearFormats()
$ws.Range("E18").Value = '  -1.08%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '11.01'
$c.ClearFormats()
$ws.Range("E19").Value = '  -0.03%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '4.27'
$c.ClearFormats()
$ws.Range("E20").Value = '  +3.37%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '328.44'
$c.ClearFormats()
$ws.Range("E21").Value = '  +0.37%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c.ClearFormats()
$ws.Range("E22").Value = '  +0.92%  '

$ws.Range("E23").Value = '  +0.01%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '10.08'
$c.ClearFormats()
$ws.Range("E24").Value = '  -1.29%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '65.75'
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '663.33'
$c.ClearFormats()
$ws.Range("E26").Value = '  +1.50%  '

$ws.Range("E27").Value = '  +5.67%  '

$ws.Range("D28").Value = '2.665.08'
$ws.Range("E28").Value = '  +3.14%  '

$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range("E29").Value = '  +0.09%  '

$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.49'
$c.ClearFormats()
$ws.Range("E30").Value = '  +2.74%  '

$ws.Range("E31").Value = '  +1.10%  '

$ws.Range("E32").Value = '  +0.97%  '

$ws.Range("E33").Value = '  +1.82%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("E35").Value = '  +1.40%  '

$ws.Range("E36").Value = '  +1.81%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.57'
$c.ClearFormats()
$ws.Range("E37").Value = '  +3.46%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '2.84'
$c.ClearFormats()
$ws.Range("E38").Value = '  +2.90%  '

$ws.Range("E39").Value = '  +0.64%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '18.92'
$c.ClearFormats()
$ws.Range("E40").Value = '  +1.07%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '151.50'
$c.ClearFormats()
$ws.Range("E41").Value = '  +0.31%  '

$ws.Range("E42").Value = '  +3.28%  '

$ws.Range("E43").Value = '  -0.01%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '159.36'
$c.ClearFormats()
$ws.Range("E44").Value = '  +3.51%  '

$ws.Range("E45").Value = '  -3.41%  '

$ws.Range("E46").Value = '  +1.41%  '

$ws.Range("E47").Value = '  +2.06%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '21.07'
$c.ClearFormats()
$ws.Range("E48").Value = '  +3.75%  '

$ws.Range("E49").Value = '  +2.32%  '

$ws.Range("E50").Value = '  +2.05%  '

$ws.Range("E51").Value = '  +2.30%  '
